# Apply updated cryptocurrency price/volume data to the symbol list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "290.07"
Set-TextValue "E2" "-1.59%"
Set-TextValue "D3" "30.74"
Set-TextValue "E3" "-4.36%"
Set-TextValue "D4" "4.915"
Set-TextValue "E4" "-2.21%"
Set-TextValue "D5" "0.07242"
Set-TextValue "E5" "-3.90%"
Set-TextValue "D6" "2.206"
Set-TextValue "E6" "15.67%"
Set-TextValue "D7" "7.648"
Set-TextValue "E7" "-2.45%"
Set-TextValue "E8" "-2.67%"
Set-TextValue "D9" "0.9000"
Set-TextValue "E9" "-2.39%"
Set-TextValue "D10" "0.1675"
Set-TextValue "E10" "-4.09%"
Set-TextValue "D11" "0.08074"
Set-TextValue "E11" "2.34%"
Set-TextValue "D12" "0.08068"
Set-TextValue "E12" "-3.91%"
Set-TextValue "D13" "0.03072"
Set-TextValue "E13" "-1.12%"
Set-TextValue "D14" "0.1001"
Set-TextValue "E14" "0.11%"
Set-TextValue "D15" "0.001505"
Set-TextValue "E15" "-0.64%"
Set-TextValue "D16" "0.005724"
Set-TextValue "E16" "-2.97%"
Set-TextValue "D17" "3.470"
Set-TextValue "E17" "0.06%"
Set-TextValue "D18" "2.075"
Set-TextValue "E18" "-3.08%"
Set-TextValue "E19" "-0.25%"
Set-TextValue "D20" "0.1303"
Set-TextValue "E20" "-1.67%"
Set-TextValue "D21" "3.968"
Set-TextValue "E21" "-10.30%"
Set-TextValue "E22" "9.49%"
Set-TextValue "E23" "-0.85%"
Set-TextValue "D24" "0.001213"
Set-TextValue "E24" "-2.20%"
Set-TextValue "D25" "0.004419"
Set-TextValue "E25" "7.89%"
Set-TextValue "D26" "0.0001299"
Set-TextValue "E26" "2.51%"
Set-TextValue "D27" "0.0003392"
Set-TextValue "E27" "-95.48%"
Set-TextValue "D39" "0.01587"
Set-TextValue "E39" "-6.81%"
Set-TextValue "D40" "0.04337"
Set-TextValue "E40" "-4.75%"
Set-TextValue "D41" "0.007280"
Set-TextValue "E41" "-0.78%"
Set-TextValue "D43" "0.1312"
Set-TextValue "E43" "-1.95%"
Set-TextValue "D44" "0.002011"
Set-TextValue "E44" "-9.45%"
Set-TextValue "D45" "0.009460"
Set-TextValue "E45" "-16.91%"
Set-TextValue "D46" "0.00005807"
Set-TextValue "E46" "-5.61%"
Set-TextValue "D47" "0.00000000749"
Set-TextValue "E47" "-1.41%"
Set-TextValue "D48" "2.255"
Set-TextValue "E48" "28.82%"
Set-TextValue "D49" "0.002897"
Set-TextValue "E49" "-4.02%"
Set-TextValue "D50" "0.00002098"
Set-TextValue "E50" "-1.41%"
Set-TextValue "D51" "0.0001998"
Set-TextValue "E51" "-1.41%"
